# Refresh the crypto tracker sheet with the latest scraped prices / 1h volumes.
#
# Column D mixes two text shapes: big coins use a dotted-thousands display
# like "42.939.70" (already non-numeric, stored as plain text) while most
# coins use an ordinary decimal like "305.22" that Excel would otherwise
# auto-convert to a Double on assignment. $TextPrefix below prepends the
# classic leading apostrophe to force those into text cells, matching the
# workbook's original (text) cell type.
$TextPrefix = "'"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.939.70'
$ws.Cells.Item(2, 5).Value = '  -0.27%  '

$ws.Cells.Item(3, 4).Value = '2.556.93'
$ws.Cells.Item(3, 5).Value = '  +0.25%  '

$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$ws.Cells.Item(5, 4).Value = $TextPrefix + '305.22'
$ws.Cells.Item(5, 5).Value = '  +2.31%  '

$ws.Cells.Item(6, 4).Value = $TextPrefix + '98.70'
$ws.Cells.Item(6, 5).Value = '  +4.75%  '

$ws.Cells.Item(7, 4).Value = $TextPrefix + '0.576'
$ws.Cells.Item(7, 5).Value = '  +0.17%  '

$ws.Cells.Item(9, 4).Value = $TextPrefix + '0.550'
$ws.Cells.Item(9, 5).Value = '  -0.12%  '

$ws.Cells.Item(10, 4).Value = $TextPrefix + '37.23'
$ws.Cells.Item(10, 5).Value = '  +2.88%  '

$ws.Cells.Item(11, 4).Value = $TextPrefix + '0.0813'
$ws.Cells.Item(11, 5).Value = '  +0.38%  '

$ws.Cells.Item(12, 2).Value = 'Polkadot'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(12, 4).Value = $TextPrefix + '7.76'
$ws.Cells.Item(12, 5).Value = '  +0.17%  '

$ws.Cells.Item(13, 2).Value = 'TRON'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(13, 4).Value = $TextPrefix + '0.116'
$ws.Cells.Item(13, 5).Value = '  +7.14%  '

$ws.Cells.Item(14, 4).Value = '2.560.68'
$ws.Cells.Item(14, 5).Value = '  +0.03%  '

$ws.Cells.Item(15, 4).Value = $TextPrefix + '15.04'
$ws.Cells.Item(15, 5).Value = '  +6.01%  '

$ws.Cells.Item(16, 4).Value = $TextPrefix + '0.886'
$ws.Cells.Item(16, 5).Value = '  +1.50%  '

$ws.Cells.Item(17, 4).Value = '42.999.63'
$ws.Cells.Item(17, 5).Value = '  -0.11%  '

$ws.Cells.Item(18, 4).Value = $TextPrefix + '13.80'
$ws.Cells.Item(18, 5).Value = '  +6.55%  '

$ws.Cells.Item(19, 5).Value = '  +0.82%  '

$ws.Cells.Item(20, 4).Value = $TextPrefix + '6.64'

$ws.Cells.Item(21, 4).Value = $TextPrefix + '71.83'
$ws.Cells.Item(21, 5).Value = '  +0.05%  '

$ws.Cells.Item(22, 4).Value = $TextPrefix + '255.06'
$ws.Cells.Item(22, 5).Value = '  -2.15%  '

$ws.Cells.Item(23, 4).Value = $TextPrefix + '2.99'
$ws.Cells.Item(23, 5).Value = '  +2.47%  '

$ws.Cells.Item(24, 5).Value = '  -1.53%  '

$ws.Cells.Item(25, 4).Value = $TextPrefix + '27.93'
$ws.Cells.Item(25, 5).Value = '  -5.51%  '

$ws.Cells.Item(26, 5).Value = '  -0.07%  '

$ws.Cells.Item(27, 4).Value = $TextPrefix + '10.17'
$ws.Cells.Item(27, 5).Value = '  +1.13%  '

$ws.Cells.Item(28, 4).Value = $TextPrefix + '38.25'
$ws.Cells.Item(28, 5).Value = '  +3.98%  '

$ws.Cells.Item(29, 4).Value = $TextPrefix + '2.19'
$ws.Cells.Item(29, 5).Value = '  +2.63%  '

$ws.Cells.Item(30, 5).Value = '  +0.63%  '

$ws.Cells.Item(31, 4).Value = $TextPrefix + '158.61'
$ws.Cells.Item(31, 5).Value = '  +2.32%  '

$ws.Cells.Item(32, 5).Value = '  +0.02%  '

$ws.Cells.Item(33, 5).Value = '  +0.84%  '

$ws.Cells.Item(34, 4).Value = $TextPrefix + '0.0811'
$ws.Cells.Item(34, 5).Value = '  +1.48%  '

$ws.Cells.Item(35, 5).Value = '  -1.55%  '

$ws.Cells.Item(36, 4).Value = $TextPrefix + '19.03'
$ws.Cells.Item(36, 5).Value = '  +15.16%  '

$ws.Cells.Item(37, 4).Value = $TextPrefix + '26.30'
$ws.Cells.Item(37, 5).Value = '  +13.51%  '

$ws.Cells.Item(38, 5).Value = '  -0.86%  '

$ws.Cells.Item(39, 5).Value = '  -0.40%  '

$ws.Cells.Item(40, 4).Value = $TextPrefix + '3.49'
$ws.Cells.Item(40, 5).Value = '  +0.42%  '

$ws.Cells.Item(41, 4).Value = $TextPrefix + '2.12'
$ws.Cells.Item(41, 5).Value = '  +33.31%  '

$ws.Cells.Item(42, 4).Value = $TextPrefix + '3.89'
$ws.Cells.Item(42, 5).Value = '  +0.02%  '

$ws.Cells.Item(43, 4).Value = '2.098.40'
$ws.Cells.Item(43, 5).Value = '  +1.41%  '

$ws.Cells.Item(44, 5).Value = '  -2.47%  '

$ws.Cells.Item(45, 4).Value = $TextPrefix + '0.998'
$ws.Cells.Item(45, 5).Value = '  +0.01%  '

$ws.Cells.Item(46, 4).Value = $TextPrefix + '86.72'
$ws.Cells.Item(46, 5).Value = '  +1.11%  '

$ws.Cells.Item(47, 4).Value = $TextPrefix + '9.11'
$ws.Cells.Item(47, 5).Value = '  +4.10%  '

$ws.Cells.Item(48, 4).Value = $TextPrefix + '75.57'
$ws.Cells.Item(48, 5).Value = '  +9.13%  '

$ws.Cells.Item(49, 4).Value = '2.803.22'
$ws.Cells.Item(49, 5).Value = '  +0.27%  '

$ws.Cells.Item(50, 4).Value = $TextPrefix + '103.86'
$ws.Cells.Item(50, 5).Value = '  -0.35%  '

$ws.Cells.Item(51, 4).Value = $TextPrefix + '0.191'
$ws.Cells.Item(51, 5).Value = '  +1.92%  '
